$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B (Step Name stays in A, a new
# "Checkout Load" step is inserted ahead of "Ultimate Load"/"Ramp Down",
# which shift right to columns C/D).
$ws.Columns("B:B").Insert()

# Carry over the existing table formatting (borders etc.) from the
# neighboring column C (previously column B) onto the newly inserted
# column B so the bordered table stays consistent.
$ws.Range("C1:C7").Copy()
$ws.Range("B1:B7").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Give the new column the same width as column A, matching the rest of
# the label/step columns.
$ws.Columns("B").ColumnWidth = $ws.Columns("A").ColumnWidth

# Header row
$ws.Range("B1").Value = "Checkout Load"
$ws.Range("C1").Value = "Ultimate Load"
$ws.Range("D1").Value = "Ramp Down"

# Row 2: Load Step Number
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 2
$ws.Range("D2").Value = 3

# Row 3: Compression [lbf]
$ws.Range("B3").Value = 800
$ws.Range("C3").Value = 2342
$ws.Range("D3").Value = 0

# Row 4: Bending [in-lbf]
$ws.Range("B4").Value = 0
$ws.Range("C4").Formula = "=976*12"
$ws.Range("D4").Value = 0

# Row 5: Bending Direction [deg]
$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 0

# Row 6: Ramp Time [s]
$ws.Range("B6").Value = 30
$ws.Range("C6").Value = 120
$ws.Range("D6").Value = 60

# Row 7: Hold Time [s]
$ws.Range("B7").Value = 10
$ws.Range("C7").Value = 60
$ws.Range("D7").Value = 30

$ws.Range("B8").Select()
